$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell $ws.Range("D2") '29.283.54'
$ws.Range("E2").Value = '  +0.55%  '

# Row 3
Set-TextCell $ws.Range("D3") '1.857.95'
$ws.Range("E3").Value = '  +0.20%  '

# Row 4
Set-TextCell $ws.Range("D4") '1.000'
$ws.Range("E4").Value = '  +0.04%  '

# Row 5
Set-TextCell $ws.Range("D5") '0.7031'
$ws.Range("E5").Value = '  +1.88%  '

# Row 6
Set-TextCell $ws.Range("D6") '238.03'
$ws.Range("E6").Value = '  +0.17%  '

# Row 7
$ws.Range("E7").Value = '  +0.11%  '

# Row 8
Set-TextCell $ws.Range("D8") '0.07938'
$ws.Range("E8").Value = '  +2.48%  '

# Row 9
Set-TextCell $ws.Range("D9") '0.3036'
$ws.Range("E9").Value = '  +0.03%  '

# Row 10
Set-TextCell $ws.Range("D10") '24.53'
$ws.Range("E10").Value = '  +7.05%  '

# Row 11
Set-TextCell $ws.Range("D11") '0.08167'
$ws.Range("E11").Value = '  +1.16%  '

# Row 12
Set-TextCell $ws.Range("D12") '1.878.89'
$ws.Range("E12").Value = '  +1.32%  '

# Row 13
Set-TextCell $ws.Range("D13") '5.219'
$ws.Range("E13").Value = '  +1.16%  '

# Row 14
Set-TextCell $ws.Range("D14") '0.7075'
$ws.Range("E14").Value = '  -1.68%  '

# Row 15
Set-TextCell $ws.Range("D15") '89.58'
$ws.Range("E15").Value = '  +0.26%  '

# Row 16
Set-TextCell $ws.Range("D16") '29.378.86'
$ws.Range("E16").Value = '  +0.86%  '

# Row 17
Set-TextCell $ws.Range("D17") '5.811'
$ws.Range("E17").Value = '  +1.61%  '

# Row 18
Set-TextCell $ws.Range("D18") '0.000007850'
$ws.Range("E18").Value = '  +1.22%  '

# Row 19
Set-TextCell $ws.Range("D19") '13.22'
$ws.Range("E19").Value = '  +0.02%  '

# Row 20
Set-TextCell $ws.Range("D20") '237.56'
$ws.Range("E20").Value = '  +1.42%  '

# Row 21
Set-TextCell $ws.Range("D21") '2.144.02'
$ws.Range("E21").Value = '  +2.48%  '

# Row 22
$ws.Range("E22").Value = '  +0.11%  '

# Row 23
Set-TextCell $ws.Range("D23") '1.000'
$ws.Range("E23").Value = '  -0.05%  '

# Row 24
Set-TextCell $ws.Range("D24") '7.572'
$ws.Range("E24").Value = '  +1.70%  '

# Row 25
Set-TextCell $ws.Range("D25") '162.57'
$ws.Range("E25").Value = '  +0.87%  '

# Row 26
Set-TextCell $ws.Range("D26") '8.914'
$ws.Range("E26").Value = '  -0.80%  '

# Row 27
$ws.Range("E27").Value = '  -0.53%  '

# Row 28
Set-TextCell $ws.Range("D28") '18.08'
$ws.Range("E28").Value = '  +0.35%  '

# Row 29
Set-TextCell $ws.Range("D29") '1.913'
$ws.Range("E29").Value = '  -1.44%  '

# Row 30
Set-TextCell $ws.Range("D30") '1.398'
$ws.Range("E30").Value = '  -0.63%  '

# Row 31
Set-TextCell $ws.Range("D31") '1.478'
$ws.Range("E31").Value = '  -0.42%  '

# Row 32
Set-TextCell $ws.Range("D32") '4.302'
$ws.Range("E32").Value = '  -3.98%  '

# Row 33
Set-TextCell $ws.Range("D33") '4.032'
$ws.Range("E33").Value = '  +0.77%  '

# Row 34
Set-TextCell $ws.Range("D34") '0.05177'
$ws.Range("E34").Value = '  -0.03%  '

# Row 35
Set-TextCell $ws.Range("D35") '1.180'
$ws.Range("E35").Value = '  +0.84%  '

# Row 36
Set-TextCell $ws.Range("D36") '0.7094'
$ws.Range("E36").Value = '  +0.40%  '

# Row 37
$ws.Range("E37").Value = '  +0.60%  '

# Row 38
$ws.Range("E38").Value = '  +0.74%  '

# Row 39
Set-TextCell $ws.Range("D39") '0.01852'
$ws.Range("E39").Value = '  +0.23%  '

# Row 40
Set-TextCell $ws.Range("D40") '2.685'

# Row 41
Set-TextCell $ws.Range("D41") '1.141.46'
$ws.Range("E41").Value = '  +3.33%  '

# Row 42
Set-TextCell $ws.Range("D42") '0.9198'
$ws.Range("E42").Value = '  -1.30%  '

# Row 43
Set-TextCell $ws.Range("D43") '5.965'
$ws.Range("E43").Value = '  +1.43%  '

# Row 44
Set-TextCell $ws.Range("D44") '0.4248'
$ws.Range("E44").Value = '  -0.45%  '

# Row 45
Set-TextCell $ws.Range("D45") '70.48'
$ws.Range("E45").Value = '  +0.14%  '

# Row 46
Set-TextCell $ws.Range("D46") '1.001'
$ws.Range("E46").Value = '  +0.05%  '

# Row 47
Set-TextCell $ws.Range("D47") '103.01'
$ws.Range("E47").Value = '  +0.41%  '

# Row 48
Set-TextCell $ws.Range("D48") '0.5309'
$ws.Range("E48").Value = '  -3.20%  '

# Row 49
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell $ws.Range("D49") '1.749'
$ws.Range("E49").Value = '  -2.03%  '

# Row 50
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell $ws.Range("D50") '9.177'
$ws.Range("E50").Value = '  +0.68%  '

# Row 51
$ws.Range("B51").Value = 'Aptos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell $ws.Range("D51") '6.995'
$ws.Range("E51").Value = '  +0.49%  '
